$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.133054
$ws.Range("H2").Value = 0.399162
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.380615
$ws.Range("N2").Value = 1.141845
$ws.Range("O2").Value = 0.04044104717146424
$ws.Range("P2").Value = 0.04044104717146424
$ws.Range("Q2").Value = 0.05064234821
$ws.Range("R2").Value = 0.45578113389
$ws.Range("S2").Value = 0.04044104717146424
$ws.Range("T2").Value = 0.04044104717146424

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.133054
$ws.Range("H3").Value = 0.399162
$ws.Range("O3").Value = 0.5170135948885495
$ws.Range("P3").Value = 0.5170135948885495
$ws.Range("Q3").Value = 0.6474308736526667
$ws.Range("R3").Value = 5.826877862874
$ws.Range("S3").Value = 0.5170135948885495
$ws.Range("T3").Value = 0.5170135948885495

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.133054
$ws.Range("H4").Value = 0.399162
$ws.Range("O4").Value = 0.4425453579399863
$ws.Range("P4").Value = 0.4425453579399863
$ws.Range("Q4").Value = 0.5541779375913333
$ws.Range("R4").Value = 4.987601438322
$ws.Range("S4").Value = 0.4425453579399863
$ws.Range("T4").Value = 0.4425453579399863
